$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: financial-period column headers (D:H) shift one year forward,
# dropping FY1396 and adding FY1401
$ws.Cells.Item(8, 4).Value  = "12 ماهه منتهی به 1397/12"
$ws.Cells.Item(8, 5).Value  = "12 ماهه منتهی به 1398/12"
$ws.Cells.Item(8, 6).Value  = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 7).Value  = "12 ماهه منتهی به 1400/12"
$ws.Cells.Item(8, 8).Value  = "12 ماهه منتهی به 1401/12"

# Row 9: matching "published on" dates (D:H)
$ws.Cells.Item(9, 4).Value  = "1399-04-30 (7)"
$ws.Cells.Item(9, 5).Value  = "1400-04-29 (9)"
$ws.Cells.Item(9, 6).Value  = "1401-04-29 (9)"
$ws.Cells.Item(9, 7).Value  = "1402-02-29 (9)"
$ws.Cells.Item(9, 8).Value  = "1402-02-29"

# Balance-sheet line items (rows 12-58, columns D:H): each period
# shifts one column left (D<-E<-F<-G<-H) and the new FY1401 figures
# are appended in column H.
$rowData = @{
    12 = @(23346, 10350, 228730, 256049, 36440)
    13 = @(0, 0, 0, 0, 0)
    14 = @(212101, 181262, 789162, 131976, 655935)
    15 = @(184520, 289886, 532034, 1076206, 2002575)
    16 = @(7585, 11498, 15546, 114998, 373369)
    17 = @(0, 0, 0, 0, 0)
    18 = @(427552, 492996, 1565472, 1579229, 3068319)
    19 = @(207, 0, 104, 0, 0)
    20 = @(64, 64, 387, 387, 387)
    21 = @(0, 0, 0, 0, 0)
    22 = @(1414611, 2959222, 5734994, 7551581, 9903974)
    23 = @(17972, 17900, 33471, 60983, 68004)
    24 = @("-", "-", "-", "-", "-")
    25 = @(2991, 120, 180, 262, 0)
    26 = @(1435845, 2977306, 5769136, 7613213, 9972365)
    27 = @(1863397, 3470302, 7334608, 9192442, 13040684)
    29 = @(445685, 992276, 978840, 995146, 2252272)
    30 = @("-", "-", "-", "-", "-")
    31 = @(6546, 40684, 258493, 601607, 519154)
    32 = @(7109, 39131, 123337, 191720, 86449)
    33 = @(2792, 100816, 6310, 222049, 352447)
    34 = @(143063, 84051, 288935, 1003498, 8074341)
    35 = @(0, 0, 0, 0, 0)
    36 = @(0, 0, 0, 0, 0)
    37 = @(605195, 1256958, 1655915, 3014020, 11284663)
    38 = @(69704, 69829, 563318, 213458, 0)
    39 = @("-", "-", "-", "-", "-")
    40 = @(702916, 1730017, 4185130, 3770292, 0)
    41 = @(4116, 4465, 6314, 7857, 15671)
    42 = @(776736, 1804311, 4754762, 3991607, 15671)
    43 = @(1381931, 3061269, 6410677, 7005627, 11300334)
    45 = @(358000, 358000, 358000, 358000, 1273000)
    46 = @(0, 0, 0, 0, 0)
    47 = @(0, 0, 0, 844638, 0)
    48 = @(0, 0, -66958, 0, 0)
    49 = @(0, 0, 0, 0, 0)
    50 = @(10924, 12672, 35800, 35800, 36692)
    51 = @(0, 0, 0, 0, 0)
    52 = @("-", "-", "-", "-", "-")
    53 = @(0, 0, 0, 0, 0)
    54 = @("-", "-", "-", "-", "-")
    55 = @(0, 0, 0, 0, 0)
    56 = @(112542, 38361, 597089, 948377, 430658)
    57 = @(481466, 409033, 923931, 2186815, 1740350)
    58 = @(1863397, 3470302, 7334608, 9192442, 13040684)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, 4 + $i).Value = $vals[$i]
    }
}
